$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")
$ws2 = $wb.Worksheets.Item("Rapport")

# --- Row 32 : fill in the end time + duration formula ---
$ws.Range("D32").Value = 0.44097222222222227
$ws.Range("E32").Formula = "=D32-C32"

# --- Row 33 : brand new entry (Conception / Documentation / Ajout des Wireframes.) ---
$ws.Range("A33").Value = 44322
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = 0.44444444444444442
$ws.Range("D33").Value = 0.46458333333333335
$ws.Range("E33").Formula = "=D33-C33"
$ws.Range("F33").Value = "Conception"
$ws.Range("G33").Value = "Documentation"
$ws.Range("H33").Value = "Ajout des Wireframes."

# --- Row 34 : new entry about SwissCenter (Communication / Échange avec le chef de projet) ---
$ws.Range("A34").Value = 44322
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = 0.46458333333333335
$ws.Range("D34").Value = 0.47569444444444442
$ws.Range("E34").Formula = "=D34-C34"
$ws.Range("F34").Value = "Communication"
$ws.Range("G34").Value = "Échange avec le chef de projet"
$ws.Range("H34").Value = "Concernant SwissCenter, le dossier de projet et le journal de travail."
$ws.Rows.Item(34).RowHeight = 30

# --- Row 35 : new entry (Conception / Documentation) ---
$ws.Range("A35").Value = 44322
$ws.Range("B35").Value = 1
$ws.Range("C35").Value = 0.47916666666666669
$ws.Range("D35").Value = 0.51041666666666663
$ws.Range("E35").Formula = "=D35-C35"
$ws.Range("F35").Value = "Conception"
$ws.Range("G35").Value = "Documentation"

# --- Row 36 : new entry, still in progress (no end time yet) ---
$ws.Range("A36").Value = 44322
$ws.Range("B36").Value = 1
$ws.Range("C36").Value = 0.52083333333333337
$ws.Range("F36").Value = "Conception"
$ws.Range("G36").Value = "Documentation"

# --- Move the current selection on the Journal sheet from H37 to H36 ---
$ws.Activate()
$ws.Range("H36").Select()

# --- Rapport sheet: scroll down a bit and move the selection to B42 ---
$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollRow = 7
$ws2.Range("B42").Select()

$ws.Activate()
